$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BecomePartnerPage")
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "BecomePartnerPage_Modal_CheckBox_CouchbaseMasterPartnerAgreement"
$ws.Range("C39").Value = "//a[contains(text(),'Couchbase Master Partner Agreement')]/parent::span/ancestor::div/input"
$ws.Activate()
$ws.Range("B39").Select()
